$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to the player's name
$ws.Name = "Robin Uthappa"

# New header row: a "matchNo" column is inserted before "teamName",
# shifting every other header one column to the right (A:L -> B:M).
$headers = @(
    "matchNo","teamName","batterName","states","runs","balls",
    "fours","sixes","sr","opponentTeamName","venue","date","result"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.ClearFormats()
}

# Full batting log - the previous single data row ("Final" vs KKR) now
# lives on row 3, and three more innings (rows 2, 4 and 5) were scraped in.
# Columns that look numeric ("runs".."sr") are written with a leading
# apostrophe so they stay text (matching the sheet's numberStoredAsText
# convention) instead of being auto-converted to real numbers.
$data = @(
    @("53rd","Chennai Super Kings","Robin Uthappa","c Harpreet Brar b Jordan","'2","'6","'0","'0","'33.33","Punjab Kings","Dubai (DSC)","October 07","Punjab Kings won by 6 wickets (with 42 balls remaining)"),
    @("Final","Chennai Super Kings","Robin Uthappa","lbw b Narine","'31","'15","'0","'3","'206.66","Kolkata Knight Riders","Dubai (DSC)","October 15","Super Kings won by 27 runs"),
    @("50th","Chennai Super Kings","Robin Uthappa","c & b Ashwin","'19","'19","'1","'0","'100.00","Delhi Capitals","Dubai (DSC)","October 04","Capitals won by 3 wickets (with 2 balls remaining)"),
    @("Qualifier","Chennai Super Kings","Robin Uthappa","c Iyer b Curran","'63","'44","'7","'2","'143.18","Delhi Capitals","Dubai (DSC)","October 10","Super Kings won by 4 wickets (with 2 balls remaining)")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        $cell.Value = $row[$c]
        $cell.ClearFormats()
    }
}
